$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.459.65'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '2.162.00'
$ws.Range('E3').Value = '  +3.16%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'227.73"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').Value = "'64.34"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.43%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('D10').Value = "'0.0859"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = "'16.05"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.47%  '
$ws.Range('D13').Value = '2.483.33'
$ws.Range('E13').Value = '  +3.22%  '
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '2.160.32'
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('D18').Value = '39.411.19'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').Value = "'71.81"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('D22').Value = "'231.65"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.56%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +5.54%  '
$ws.Range('D25').Value = "'2.36"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').Value = "'172.51"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').Value = "'19.88"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('D31').Value = "'2.67"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.39%  '
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').Value = "'4.77"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  +8.83%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').Value = "'3.57"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').Value = "'1.00"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').Value = "'104.01"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = "'17.90"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D43').Value = '1.538.92'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = "'1.18"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.71%  '
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('E47').Value = '  +5.74%  '
$ws.Range('D48').Value = "'7.79"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = "'4.19"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').Value = '2.365.61'
$ws.Range('E50').Value = '  +3.30%  '
$ws.Range('E51').Value = '  -0.03%  '
